# Iraq League base update - 28-05-2024 19:13
# 1) Eight pairs of adjacent match rows had their B,C,E:AD data swapped
#    (column A - the running id - and column D - the match date, which is
#    identical within each pair - are left untouched).
# 2) Three brand-new match rows (223-225) were appended at the bottom of
#    the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($rowA, $rowB) {
    $bcA = $ws.Range("B$($rowA):C$($rowA)").Value()
    $edA = $ws.Range("E$($rowA):AD$($rowA)").Value()
    $bcB = $ws.Range("B$($rowB):C$($rowB)").Value()
    $edB = $ws.Range("E$($rowB):AD$($rowB)").Value()

    $ws.Range("B$($rowA):C$($rowA)").Value = $bcB
    $ws.Range("E$($rowA):AD$($rowA)").Value = $edB

    $ws.Range("B$($rowB):C$($rowB)").Value = $bcA
    $ws.Range("E$($rowB):AD$($rowB)").Value = $edA
}

Swap-MatchRows 17 18
Swap-MatchRows 41 42
Swap-MatchRows 58 59
Swap-MatchRows 69 70
Swap-MatchRows 73 74
Swap-MatchRows 135 136
Swap-MatchRows 151 152
Swap-MatchRows 173 174

# Append the three new fixtures below the previous last row (222).

$ws.Range("A222:C222").Copy()
$ws.Range("A223:C225").PasteSpecial(-4122)
$ws.Range("D222").Copy()
$ws.Range("D223:D225").PasteSpecial(-4122)
$ws.Range("E222:H222").Copy()
$ws.Range("E223:H225").PasteSpecial(-4122)
$ws.Range("K222:AD222").Copy()
$ws.Range("K223:AD225").PasteSpecial(-4122)

$ws.Range("A223").Value = 221
$ws.Range("B223").Value = 8267113
$ws.Range("C223").Value = "Iraq League"
$ws.Range("D223").Value = 45439.5
$ws.Range("E223").Value = "Karbalaa FC"
$ws.Range("F223").Value = "Al Najaf"
$ws.Range("G223").Value = 0
$ws.Range("H223").Value = 0
$ws.Range("K223").Value = "D"
$ws.Range("L223").Value = 3.75
$ws.Range("M223").Value = 2.875
$ws.Range("N223").Value = 2
$ws.Range("O223").Value = 3.2
$ws.Range("P223").Value = 2.8
$ws.Range("Q223").Value = 2.25
$ws.Range("R223").Value = 0.25
$ws.Range("S223").Value = 1.825
$ws.Range("T223").Value = 1.975
$ws.Range("U223").Value = 2
$ws.Range("V223").Value = 1.975
$ws.Range("W223").Value = 1.825
$ws.Range("X223").Value = -1
$ws.Range("Y223").Value = 1.8
$ws.Range("Z223").Value = -1
$ws.Range("AA223").Value = 0.4125
$ws.Range("AB223").Value = -0.5
$ws.Range("AC223").Value = -1
$ws.Range("AD223").Value = 0.825

$ws.Range("A224").Value = 222
$ws.Range("B224").Value = 8263267
$ws.Range("C224").Value = "Iraq League"
$ws.Range("D224").Value = 45439.5
$ws.Range("E224").Value = "Zakho"
$ws.Range("F224").Value = "Al Naft SC"
$ws.Range("G224").Value = 3
$ws.Range("H224").Value = 0
$ws.Range("K224").Value = "H"
$ws.Range("L224").Value = 2.2
$ws.Range("M224").Value = 2.8
$ws.Range("N224").Value = 3.25
$ws.Range("O224").Value = 2.375
$ws.Range("P224").Value = 2.625
$ws.Range("Q224").Value = 3.1
$ws.Range("R224").Value = -0.25
$ws.Range("S224").Value = 2.025
$ws.Range("T224").Value = 1.775
$ws.Range("U224").Value = 2
$ws.Range("V224").Value = 2.025
$ws.Range("W224").Value = 1.775
$ws.Range("X224").Value = 1.375
$ws.Range("Y224").Value = -1
$ws.Range("Z224").Value = -1
$ws.Range("AA224").Value = 1.025
$ws.Range("AB224").Value = -1
$ws.Range("AC224").Value = 1.025
$ws.Range("AD224").Value = -1

$ws.Range("A225").Value = 223
$ws.Range("B225").Value = 8263268
$ws.Range("C225").Value = "Iraq League"
$ws.Range("D225").Value = 45439.58333333334
$ws.Range("E225").Value = "Al Quwa Al Jawiya"
$ws.Range("F225").Value = "Al Zawraa"
$ws.Range("G225").Value = 0
$ws.Range("H225").Value = 1
$ws.Range("K225").Value = "A"
$ws.Range("L225").Value = 2.25
$ws.Range("M225").Value = 2.75
$ws.Range("N225").Value = 3.25
$ws.Range("O225").Value = 2.05
$ws.Range("P225").Value = 2.75
$ws.Range("Q225").Value = 3.7
$ws.Range("R225").Value = -0.25
$ws.Range("S225").Value = 1.775
$ws.Range("T225").Value = 2.025
$ws.Range("U225").Value = 2
$ws.Range("V225").Value = 1.925
$ws.Range("W225").Value = 1.875
$ws.Range("X225").Value = -1
$ws.Range("Y225").Value = -1
$ws.Range("Z225").Value = 2.7
$ws.Range("AA225").Value = -1
$ws.Range("AB225").Value = 1.025
$ws.Range("AC225").Value = -1
$ws.Range("AD225").Value = 0.875

# I / J (half-time score) columns were never recorded for these three
# fixtures, mirroring row 222 immediately above them.
$ws.Range("I223:J225").ClearContents()
